# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
# Mirrors the commit: "Created functions to get season record" — the
# per-player rows each get the team's overall W/L/T totals for the season
# appended after the existing "Unnamed: 28" column (AC), extending the
# used range from A1:AC63 to A1:AF63.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastDataRow = 63
$winsCol   = 30   # AD
$lossesCol = 31   # AE
$tiesCol   = 32   # AF

# Give the three new header cells the same look (bold / bordered / centered)
# as the rest of row 1 by copying the formatting from the neighboring
# header cell (AC1) before filling in the text.
$ws.Cells.Item(1, 29).Copy()
$headerRange = $ws.Range($ws.Cells.Item(1, $winsCol), $ws.Cells.Item(1, $tiesCol))
$headerRange.PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(1, $winsCol).Value   = "Wins"
$ws.Cells.Item(1, $lossesCol).Value = "Losses"
$ws.Cells.Item(1, $tiesCol).Value   = "Ties"

# Every player row gets the same season record for the team.
for ($r = 2; $r -le $lastDataRow; $r++) {
    $ws.Cells.Item($r, $winsCol).Value   = 101
    $ws.Cells.Item($r, $lossesCol).Value = 61
    $ws.Cells.Item($r, $tiesCol).Value   = 0
}
